$wb = $excel.ActiveWorkbook

$wsPortfolio = $wb.Worksheets.Item("Portfolio Design")
$wsFeedstock = $wb.Worksheets.Item("feedstock_to_commodity")

# Rebalance the "Diverted Organic Waste" vs "Animal Manure" biomethane split
# (6 MCF total kept the same, but shifted almost entirely to Animal Manure).
$wsPortfolio.Range("C28").Value = 5.7692307692307692
$wsPortfolio.Range("C29").Value = 0.23076923076923078

# Mirror the author's on-screen review: glance at the recalculated
# feedstock_to_commodity totals (Z15) before returning to the edited cells.
$wsFeedstock.Activate()
$wsFeedstock.Range("Z15").Select() | Out-Null

$wsPortfolio.Activate()
$wsPortfolio.Range("C28:C29").Select() | Out-Null
